$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the two new worksheets at the end of the workbook, in order:
#   fabric-ap-int-policy (sheet16), fabric-ap-int-block (sheet17)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPolicy = $wb.Worksheets.Add($null, $lastSheet)
$wsPolicy.Name = "fabric-ap-int-policy"

$wsBlock = $wb.Worksheets.Add($null, $wsPolicy)
$wsBlock.Name = "fabric-ap-int-block"

# ---------------------------------------------------------------------------
# fabric-ap-int-policy data
# ---------------------------------------------------------------------------
$policyData = @(
    @("intpolicy", "linktype", "speedpolicy", "cdppolicy", "lldppolicy", "pcpolicy", "state"),
    @("Firewall-policy", "leaf", "Speed_1G", "default", "default", "", "present"),
    @("LB-policy", "leaf", "Speed_1G", "default", "default", "", "present"),
    @("c7000-policy", "node", "Speed_1G", "default", "default", "lacp_active", "present"),
    @("Chassis-policy", "link", "Speed_1G", "default", "default", "default", "present")
)

for ($r = 0; $r -lt $policyData.Count; $r++) {
    $row = $policyData[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $val = $row[$c]
        if ($val -ne "") {
            $wsPolicy.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# ---------------------------------------------------------------------------
# fabric-ap-int-block data
# ---------------------------------------------------------------------------
$blockHeaders = @("intleafprf", "intname", "intfrom", "intto", "intpolicy", "state")
for ($c = 0; $c -lt $blockHeaders.Count; $c++) {
    $wsBlock.Cells.Item(1, $c + 1).Value = $blockHeaders[$c]
}

$blockData = @(
    @("INT-PRF-LF-1101", "Eth1_1", 1, 1, "Firewall-policy", "present"),
    @("INT-PRF-LF-1102", "Eth1_10", 10, 10, "Firewall-policy", "present"),
    @("INT-PRF-LF-1103", "Eth1_3", 3, 3, "LB-policy", "present"),
    @("INT-PRF-LF-1104", "Eth1_5", 5, 5, "LB-policy", "present"),
    @("INT-PRF-LF-1101-1102", "Eth1_2", 2, 2, "c7000-policy", "present"),
    @("INT-PRF-LF-1101-1102", "Eth1_22", 22, 22, "c7000-policy", "present"),
    @("INT-PRF-LF-1103-1104", "Eth1_9", 9, 9, "Chassis-policy", "present"),
    @("INT-PRF-LF-1103-1104", "Eth1_8", 8, 8, "Chassis-policy", "present")
)

for ($r = 0; $r -lt $blockData.Count; $r++) {
    $row = $blockData[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsBlock.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------------
# Column widths (approximate Excel "character width" units used by the
# original workbook) for the two new sheets.
# ---------------------------------------------------------------------------
$policyWidths = @(12.96, 7.54, 7.54, 7.54, 7.54, 7.54, 7.54)
for ($c = 0; $c -lt $policyWidths.Count; $c++) {
    $wsPolicy.Columns.Item($c + 1).ColumnWidth = $policyWidths[$c]
}

$blockWidths = @(19.77, 7.95, 7.13, 5.04, 12.96, 7.54)
for ($c = 0; $c -lt $blockWidths.Count; $c++) {
    $wsBlock.Columns.Item($c + 1).ColumnWidth = $blockWidths[$c]
}

# ---------------------------------------------------------------------------
# Selection / active-cell bookkeeping to mirror the authored workbook: the
# policy sheet was left mid-edit at D12, the block sheet (last-active tab)
# sits at A1.
# ---------------------------------------------------------------------------
$null = $wsPolicy.Range("D12").Select()
$null = $wsBlock.Range("A1").Select()
$null = $wsBlock.Activate()
